$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("F18").Value = "['ELM-1NA-Tecnologia dos Materiais.', 'ELM-1NA-Tecnologia dos Materiais.']"

# Row 19
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "['MEC-1NB-T.M. Metalicos', -, -, -]"
$ws.Range("E19").Value = "[-, 'MEC-1NA-Trat. Termicos', -, -]"

# Row 20
$ws.Range("B20").Value = "[-, -, -, 'MEC-1NA-Trat. Termicos']"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "['MEC-1NB-T.M. Metalicos', -, -, -]"
$ws.Range("E20").Value = "[-, 'MEC-1NA-Trat. Termicos', -, -]"
$ws.Range("F20").Value = "[-, -, -, 'MEC-1NB-T.M. Metalicos']"

# Row 21
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "[-, 'MEC-1NA-Trat. Termicos', -, -]"
$ws.Range("F21").Value = "['MEC-1NB-T.M. Metalicos', -, -, -]"
